# Atualizacao de bases das ligas (Mexico Liga de Expansion), 17-02-2024 11:11
# - Rows 47/48, 102/103 and 262/263 each had their data (all columns except id
#   in column A) swapped between the two rows.
# - Row 334 is overwritten with what used to be row 335's data, row 335 is
#   overwritten with updated odds for the match that used to live in row 336,
#   and the now-redundant row 336 is deleted (rows shift up, dimension shrinks
#   from AC336 to AC335).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 47
$ws.Range("B47").Value = 6007835
$ws.Range("F47").Value = "Cimarrones de Sonora FC"
$ws.Range("G47").Value = "Tlaxcala FC"
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 1
$ws.Range("J47").Value = "A"
$ws.Range("K47").Value = 1.7
$ws.Range("L47").Value = 3.5
$ws.Range("M47").Value = 4.5
$ws.Range("N47").Value = 1.363
$ws.Range("O47").Value = 4.75
$ws.Range("P47").Value = 8.5
$ws.Range("Q47").Value = -1.5
$ws.Range("R47").Value = 1.95
$ws.Range("S47").Value = 1.85
$ws.Range("T47").Value = 2.75
$ws.Range("U47").Value = 1.825
$ws.Range("V47").Value = 1.975
$ws.Range("W47").Value = -1
$ws.Range("X47").Value = -1
$ws.Range("Y47").Value = 7.5
$ws.Range("Z47").Value = -1
$ws.Range("AA47").Value = 0.8500000000000001
$ws.Range("AB47").Value = -1
$ws.Range("AC47").Value = 0.9750000000000001

# Row 48
$ws.Range("B48").Value = 6007834
$ws.Range("F48").Value = "Universidad Guadalajara"
$ws.Range("G48").Value = "Alacranes de Durango"
$ws.Range("H48").Value = 1
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = "H"
$ws.Range("K48").Value = 1.909
$ws.Range("L48").Value = 3.3
$ws.Range("M48").Value = 3.75
$ws.Range("N48").Value = 1.909
$ws.Range("O48").Value = 3.4
$ws.Range("P48").Value = 4
$ws.Range("Q48").Value = -0.5
$ws.Range("R48").Value = 1.925
$ws.Range("S48").Value = 1.875
$ws.Range("T48").Value = 2.25
$ws.Range("U48").Value = 1.85
$ws.Range("V48").Value = 1.95
$ws.Range("W48").Value = 0.909
$ws.Range("X48").Value = -1
$ws.Range("Y48").Value = -1
$ws.Range("Z48").Value = 0.925
$ws.Range("AA48").Value = -1
$ws.Range("AB48").Value = -1
$ws.Range("AC48").Value = 0.95

# Row 102
$ws.Range("B102").Value = 6007882
$ws.Range("F102").Value = "Cimarrones de Sonora FC"
$ws.Range("G102").Value = "Club Celaya"
$ws.Range("H102").Value = 1
$ws.Range("I102").Value = 1
$ws.Range("J102").Value = "D"
$ws.Range("K102").Value = 2.1
$ws.Range("L102").Value = 3.2
$ws.Range("M102").Value = 3.25
$ws.Range("N102").Value = 2.2
$ws.Range("O102").Value = 3
$ws.Range("P102").Value = 3.6
$ws.Range("Q102").Value = -0.25
$ws.Range("R102").Value = 1.875
$ws.Range("S102").Value = 1.925
$ws.Range("T102").Value = 2
$ws.Range("U102").Value = 1.8
$ws.Range("V102").Value = 2
$ws.Range("W102").Value = -1
$ws.Range("X102").Value = 2
$ws.Range("Y102").Value = -1
$ws.Range("Z102").Value = -0.5
$ws.Range("AA102").Value = 0.4625
$ws.Range("AB102").Value = 0
$ws.Range("AC102").Value = 0

# Row 103
$ws.Range("B103").Value = 6007883
$ws.Range("F103").Value = "Club Atletico La Paz"
$ws.Range("G103").Value = "Dorados"
$ws.Range("H103").Value = 1
$ws.Range("I103").Value = 1
$ws.Range("J103").Value = "D"
$ws.Range("K103").Value = 1.65
$ws.Range("L103").Value = 3.75
$ws.Range("M103").Value = 4.5
$ws.Range("N103").Value = 1.533
$ws.Range("O103").Value = 4.333
$ws.Range("P103").Value = 6
$ws.Range("Q103").Value = -1
$ws.Range("R103").Value = 1.825
$ws.Range("S103").Value = 1.975
$ws.Range("T103").Value = 2.75
$ws.Range("U103").Value = 1.95
$ws.Range("V103").Value = 1.85
$ws.Range("W103").Value = -1
$ws.Range("X103").Value = 3.333
$ws.Range("Y103").Value = -1
$ws.Range("Z103").Value = -1
$ws.Range("AA103").Value = 0.9750000000000001
$ws.Range("AB103").Value = -1
$ws.Range("AC103").Value = 0.8500000000000001

# Row 262
$ws.Range("B262").Value = 6924569
$ws.Range("F262").Value = "Venados FC"
$ws.Range("G262").Value = "Dorados"
$ws.Range("H262").Value = 4
$ws.Range("I262").Value = 1
$ws.Range("J262").Value = "H"
$ws.Range("K262").Value = 1.615
$ws.Range("L262").Value = 4
$ws.Range("M262").Value = 4.5
$ws.Range("N262").Value = 1.5
$ws.Range("O262").Value = 4.75
$ws.Range("P262").Value = 5.75
$ws.Range("Q262").Value = -1.25
$ws.Range("R262").Value = 1.925
$ws.Range("S262").Value = 1.875
$ws.Range("T262").Value = 3
$ws.Range("U262").Value = 1.75
$ws.Range("V262").Value = 1.95
$ws.Range("W262").Value = 0.5
$ws.Range("X262").Value = -1
$ws.Range("Y262").Value = -1
$ws.Range("Z262").Value = 0.925
$ws.Range("AA262").Value = -1
$ws.Range("AB262").Value = 0.75
$ws.Range("AC262").Value = -1

# Row 263
$ws.Range("B263").Value = 6924568
$ws.Range("F263").Value = "Atletico Morelia"
$ws.Range("G263").Value = "Atlante"
$ws.Range("H263").Value = 0
$ws.Range("I263").Value = 1
$ws.Range("J263").Value = "A"
$ws.Range("K263").Value = 2.4
$ws.Range("L263").Value = 3
$ws.Range("M263").Value = 2.875
$ws.Range("N263").Value = 2.7
$ws.Range("O263").Value = 3.1
$ws.Range("P263").Value = 2.8
$ws.Range("Q263").Value = 0
$ws.Range("R263").Value = 1.85
$ws.Range("S263").Value = 1.95
$ws.Range("T263").Value = 2.25
$ws.Range("U263").Value = 1.975
$ws.Range("V263").Value = 1.725
$ws.Range("W263").Value = -1
$ws.Range("X263").Value = -1
$ws.Range("Y263").Value = 1.8
$ws.Range("Z263").Value = -1
$ws.Range("AA263").Value = 0.95
$ws.Range("AB263").Value = -1
$ws.Range("AC263").Value = 0.7250000000000001

# Row 334
$ws.Range("B334").Value = 7641673
$ws.Range("E334").Value = 45339.83680555555
$ws.Range("F334").Value = "Tapatio"
$ws.Range("G334").Value = "Atletico Morelia"
$ws.Range("K334").Value = 2.45
$ws.Range("L334").Value = 3.4
$ws.Range("M334").Value = 2.6
$ws.Range("N334").Value = 2.4
$ws.Range("O334").Value = 3.3
$ws.Range("P334").Value = 3
$ws.Range("Q334").Value = -0.25
$ws.Range("R334").Value = 2.05
$ws.Range("S334").Value = 1.75
$ws.Range("T334").Value = 2.5
$ws.Range("U334").Value = 2.025
$ws.Range("V334").Value = 1.775
$ws.Range("W334").Value = 0
$ws.Range("X334").Value = 0
$ws.Range("Y334").Value = 0
$ws.Range("Z334").Value = 0
$ws.Range("AA334").Value = 0

# Row 335
$ws.Range("B335").Value = 7640645
$ws.Range("E335").Value = 45340.625
$ws.Range("F335").Value = "Atlante"
$ws.Range("G335").Value = "Cimarrones de Sonora FC"
$ws.Range("K335").Value = 1.55
$ws.Range("L335").Value = 3.75
$ws.Range("M335").Value = 5.75
$ws.Range("N335").Value = 1.571
$ws.Range("O335").Value = 4
$ws.Range("P335").Value = 6
$ws.Range("Q335").Value = -1
$ws.Range("R335").Value = 1.9
$ws.Range("S335").Value = 1.9
$ws.Range("T335").Value = 2.5
$ws.Range("U335").Value = 1.95
$ws.Range("V335").Value = 1.85
$ws.Range("W335").Value = 0
$ws.Range("X335").Value = 0
$ws.Range("Y335").Value = 0
$ws.Range("Z335").Value = 0
$ws.Range("AA335").Value = 0

# Delete row 336 (data merged into row 335, row shifted out)
$ws.Rows.Item(336).Delete()

$wb.Save()